# Update the "tblresumenfactura" summary row (row 2) with the new invoice
# totals. Most of these values are stored as literal text (matching the
# exporter's original formatting, e.g. trailing zeros like "115267.32000"),
# so the target ranges are pre-formatted as Text before the value is typed
# in - this stops Excel from auto-converting the numeric-looking strings
# into real numbers. TotalIVADevuelto (M2) is the one column that becomes
# an actual number (0) instead of text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells A2:L2 and N2:P2 hold text values - mark them as Text format first
# so the values we assign stick as strings instead of being reinterpreted
# as numbers.
$ws.Range("A2:L2").NumberFormat = "@"
$ws.Range("N2:P2").NumberFormat = "@"

$ws.Range("A2").Value = "0"
$ws.Range("B2").Value = "0"
$ws.Range("C2").Value = "0"
$ws.Range("D2").Value = "115267.32000"
$ws.Range("E2").Value = "0"
$ws.Range("F2").Value = "0"
$ws.Range("G2").Value = "115267.32000"
$ws.Range("H2").Value = "0"
$ws.Range("I2").Value = "0"
$ws.Range("J2").Value = "115267.32000"
$ws.Range("K2").Value = "102152.99000"
$ws.Range("L2").Value = "13761.02000"
$ws.Range("M2").Value = 0
$ws.Range("N2").Value = "115914.01000"
$ws.Range("O2").Value = "0"
$ws.Range("P2").Value = "13114.33000"

# Restore the default (General) style on the text cells now that the
# literal text values are locked in, so no visible formatting change is
# left behind on the worksheet.
$ws.Range("A2:L2").Style = "Normal"
$ws.Range("N2:P2").Style = "Normal"
